$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from the very last paragraph ("Third item"
#    in "Another unordered list") to just after "Second item" in the
#    "Unordered list" section (paragraph 14 in the original document).
#
#    Word places a collapsed "_GoBack" bookmark at the position of the most
#    recent edit. Adding a bookmark with the same name elsewhere moves it
#    (bookmark names are unique), so the stale one at the end is removed
#    automatically.
#
#    Placing a *collapsed* bookmark exactly at "end of paragraph text" (i.e.
#    immediately before the paragraph mark) is done indirectly: insert a
#    throwaway character there, bookmark the position that is now in front
#    of it (no longer the paragraph's last position), then delete the
#    throwaway character again. The bookmark survives at the right spot.
# ---------------------------------------------------------------------------

$targetPara = $d.Paragraphs(14)
$targetEnd = $targetPara.Range.End - 1   # position right before the pilcrow

$marker = $d.Range($targetEnd, $targetEnd)
$marker.InsertAfter("X")

$bookmarkSpot = $d.Range($targetEnd, $targetEnd)
$d.Bookmarks.Add("_GoBack", $bookmarkSpot)

$markerChar = $d.Range($targetEnd, $targetEnd + 1)
$markerChar.Delete()

# ---------------------------------------------------------------------------
# 2) Remove the empty paragraph between the "Second item" paragraph above and
#    the "Another ordered list" heading (originally paragraph 15).
# ---------------------------------------------------------------------------
$d.Paragraphs(15).Range.Delete()

# ---------------------------------------------------------------------------
# 3) Remove the empty paragraph between the first "Third item" (end of the
#    "Ordered list" section) and the "Unordered list" heading (originally
#    paragraph 10). Done last since it sits before the other edits and would
#    otherwise shift their indices.
# ---------------------------------------------------------------------------
$d.Paragraphs(10).Range.Delete()
